# chore: update Sheets via scheduled runner
# Refresh cached market-board price/profit figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# leve-profit sheets. A handful of rows also lose a stale profit/loss
# column (cell removed entirely) where the corresponding price is now 0.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 189
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

$ws.Range("H12").Value = 574.25
$ws.Range("I12").Value = 299
$ws.Range("J12").Value = 666
$ws.Range("K12").Value = 299
$ws.Range("L12").Value = 666
$ws.Range("M12").Value = -129
$ws.Range("N12").Value = -1006

$ws.Range("H19").Value = 1079.6666
$ws.Range("I19").Value = 700
$ws.Range("K19").Value = 700
$ws.Range("M19").Value = -525

$ws.Range("H103").Value = 674.53845
$ws.Range("I103").Value = 583
$ws.Range("J103").Value = 781.3333
$ws.Range("K103").Value = 1749
$ws.Range("L103").Value = 2343.9999
$ws.Range("M103").Value = -1163
$ws.Range("N103").Value = -3515.9999

$ws.Range("H135").Value = 20835188
$ws.Range("I135").Value = 21740936
$ws.Range("K135").Value = 195668424
$ws.Range("M135").Value = -195665889

$ws.Range("H138").Value = 3449.4614
$ws.Range("I138").Value = 2534.2144
$ws.Range("K138").Value = 7602.6432
$ws.Range("M138").Value = -2462.6432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 19525.389
$ws.Range("I2").Value = 30232.637
$ws.Range("J2").Value = 2699.7144
$ws.Range("K2").Value = 30232.637
$ws.Range("L2").Value = 2699.7144
$ws.Range("M2").Value = -30119.637
$ws.Range("N2").Value = -2925.7144

$ws.Range("H45").Value = 5001.385
$ws.Range("I45").Value = 3252.375
$ws.Range("J45").Value = 7799.8
$ws.Range("K45").Value = 3252.375
$ws.Range("L45").Value = 7799.8
$ws.Range("M45").Value = -2875.375
$ws.Range("N45").Value = -8553.799999999999

$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H61").Value = 3000
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

$ws.Range("H74").Value = 1736.4722
$ws.Range("I74").Value = 1795.72
$ws.Range("K74").Value = 1795.72
$ws.Range("M74").Value = -921.72

$ws.Range("H77").Value = 1736.4722
$ws.Range("I77").Value = 1795.72
$ws.Range("K77").Value = 8978.6
$ws.Range("M77").Value = -4610.6

$ws.Range("H116").Value = 19525.389
$ws.Range("I116").Value = 30232.637
$ws.Range("J116").Value = 2699.7144
$ws.Range("K116").Value = 30232.637
$ws.Range("L116").Value = 2699.7144
$ws.Range("M116").Value = -27938.637
$ws.Range("N116").Value = -7287.7144

$ws.Range("H132").Value = 2571.0454
$ws.Range("I132").Value = 2517.5557
$ws.Range("J132").Value = 2811.75
$ws.Range("K132").Value = 7552.6671
$ws.Range("L132").Value = 8435.25
$ws.Range("M132").Value = -5022.6671
$ws.Range("N132").Value = -13495.25

$ws.Range("H136").Value = 3000
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 19525.389
$ws.Range("I3").Value = 30232.637
$ws.Range("J3").Value = 2699.7144
$ws.Range("K3").Value = 30232.637
$ws.Range("L3").Value = 2699.7144
$ws.Range("M3").Value = -30118.637
$ws.Range("N3").Value = -2927.7144

$ws.Range("H12").Value = 224.66667
$ws.Range("I12").Value = 69.8
$ws.Range("J12").Value = 999
$ws.Range("K12").Value = 69.8
$ws.Range("L12").Value = 999
$ws.Range("M12").Value = 98.2
$ws.Range("N12").Value = -1335

$ws.Range("H20").Value = 2862.4285
$ws.Range("I20").Value = 2126.8572
$ws.Range("J20").Value = 3598
$ws.Range("K20").Value = 2126.8572
$ws.Range("L20").Value = 3598
$ws.Range("M20").Value = -1879.8572
$ws.Range("N20").Value = -4092

$ws.Range("H86").Value = 17243376
$ws.Range("I86").Value = 20002098
$ws.Range("J86").Value = 1363.75
$ws.Range("K86").Value = 20002098
$ws.Range("L86").Value = 1363.75
$ws.Range("M86").Value = -20000975
$ws.Range("N86").Value = -3609.75

$ws.Range("H89").Value = 17243376
$ws.Range("I89").Value = 20002098
$ws.Range("J89").Value = 1363.75
$ws.Range("K89").Value = 100010490
$ws.Range("L89").Value = 6818.75
$ws.Range("M89").Value = -100004874
$ws.Range("N89").Value = -18050.75

$ws.Range("H100").Value = 14812.25
$ws.Range("J100").Value = 14812.25
$ws.Range("L100").Value = 14812.25
$ws.Range("N100").Value = -16976.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 14695.4
$ws.Range("I38").Value = 14695.4
$ws.Range("K38").Value = 14695.4
$ws.Range("M38").Value = -14318.4

$ws.Range("H46").Value = 14695.4
$ws.Range("I46").Value = 14695.4
$ws.Range("K46").Value = 14695.4
$ws.Range("M46").Value = -14484.4

$ws.Range("H62").Value = 62502450
$ws.Range("I62").Value = 3799.75
$ws.Range("J62").Value = 125001096
$ws.Range("K62").Value = 3799.75
$ws.Range("L62").Value = 125001096
$ws.Range("M62").Value = -3175.75
$ws.Range("N62").Value = -125002344

$ws.Range("H65").Value = 62502450
$ws.Range("I65").Value = 3799.75
$ws.Range("J65").Value = 125001096
$ws.Range("K65").Value = 18998.75
$ws.Range("L65").Value = 625005480
$ws.Range("M65").Value = -15878.75
$ws.Range("N65").Value = -625011720

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 629.8
$ws.Range("I9").Value = 666.3333
$ws.Range("K9").Value = 1998.9999
$ws.Range("M9").Value = -1774.9999

$ws.Range("H137").Value = 2458.3333
$ws.Range("J137").Value = 2416.6667
$ws.Range("L137").Value = 7250.000100000001
$ws.Range("N137").Value = -17450.0001

$ws.Range("H139").Value = 125002710
$ws.Range("I139").Value = 125002710
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 375008130
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -375002990
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 19363.5
$ws.Range("J96").Value = 19363.5
$ws.Range("L96").Value = 19363.5
$ws.Range("N96").Value = -24855.5

$ws.Range("H102").Value = 3734.8928
$ws.Range("I102").Value = 3423.1904
$ws.Range("K102").Value = 3423.1904
$ws.Range("M102").Value = -1801.1904

$ws.Range("H126").Value = 3259.2
$ws.Range("I126").Value = 3756
$ws.Range("J126").Value = 2928
$ws.Range("K126").Value = 11268
$ws.Range("L126").Value = 8784
$ws.Range("M126").Value = -8798
$ws.Range("N126").Value = -13724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3154.4
$ws.Range("I46").Value = 2123.75
$ws.Range("J46").Value = 3841.5
$ws.Range("K46").Value = 2123.75
$ws.Range("L46").Value = 3841.5
$ws.Range("M46").Value = -1935.75
$ws.Range("N46").Value = -4217.5

$ws.Range("H106").Value = 8465.588
$ws.Range("J106").Value = 8465.588
$ws.Range("L106").Value = 8465.588
$ws.Range("N106").Value = -10989.588

$ws.Range("H122").Value = 3705.5789
$ws.Range("I122").Value = 3067.8965
$ws.Range("K122").Value = 9203.6895
$ws.Range("M122").Value = -6753.6895

$ws.Range("H132").Value = 2302.4614
$ws.Range("I132").Value = 1907.2273
$ws.Range("K132").Value = 5721.6819
$ws.Range("M132").Value = -3191.6819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 13131.444
$ws.Range("I45").Value = 4468
$ws.Range("J45").Value = 15606.714
$ws.Range("K45").Value = 4468
$ws.Range("L45").Value = 15606.714
$ws.Range("M45").Value = -3977
$ws.Range("N45").Value = -16588.714

$ws.Range("H132").Value = 1326.2273
$ws.Range("J132").Value = 6000
$ws.Range("L132").Value = 18000
$ws.Range("N132").Value = -23060
